$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split-payment test data: custcredit value (column X, row 2) changes from 100 to 0.
# Leading apostrophe keeps it text-typed (matches the original quote-prefixed cell style).
$ws.Range("X2").Value = "'0"

# Selection moved from N5 to X3 (and the view scrolled so column R is now left-most).
$ws.Range("X3").Select()
